# Auto-update draw results: append the 2025-10-24 "Pick 4" row.
#
# The sheet stores every value as literal text (dates like "2025-10-24" and
# phase codes like "251024" are plain strings, not numbers/dates). Excel's
# COM layer auto-detects "number-looking" and "date-looking" text typed into
# a cell via .Value and silently coerces it to a real number/date. Prefixing
# with an apostrophe forces Excel to keep it as literal text (the apostrophe
# itself is not stored in the value) for the two columns where that matters
# (A = date-like, C = numeric-like). Columns B/D/E are not number/date-like
# so they remain plain text without needing the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = "'2025-10-24"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "'251024"
$ws.Cells.Item($row, 4).Value = "2-5-0-1"
$ws.Cells.Item($row, 5).Value = "2025-10-24T21:38:03.959+04:00"
